# Wait until movement complete Added
#
# The original "Sheet1" (the per-run test template/results sheet) is
# renamed to "Tested" (it already holds the results of a completed run:
# video_file / gz_pose_file / vid_pose_file flags). A fresh copy of that
# sheet is then added, named "Sheet1" again, acting as the blank template
# for the next test run (all three flags reset to 0). The new "Sheet1"
# becomes the active tab.

$wb = $excel.ActiveWorkbook

$orig = $wb.Worksheets.Item("Sheet1")

# Update the completed-run flags (video_file / gz_pose_file / vid_pose_file)
# from their old free-text placeholders to numeric status flags.
$orig.Range("B3").Value = 0
$orig.Range("B4").Value = 1
$orig.Range("B5").Value = 0
$null = $orig.Range("B4").Select()

# Rename the completed run sheet to "Tested".
$orig.Name = "Tested"

# Duplicate it to create the fresh template for the next run, placed right
# after "Tested".
$orig.Copy($null, $orig)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Sheet1"

# Reset the new template's status flags back to "not yet run".
$newSheet.Range("B3").Value = 0
$newSheet.Range("B4").Value = 0
$newSheet.Range("B5").Value = 0

# New template sheet is the active / selected sheet.
$null = $newSheet.Range("B5").Select()
